# Updated cryptos list on Wed Nov  1 16:33:31 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "34.522.19"
$ws.Range("E2").Value = "  +0.54%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.802.38"
$ws.Range("E3").Value = "  +0.31%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.25%  "

# Row 5 - BNB
$ws.Range("D5").Value = "224.34"
$ws.Range("E5").Value = "  -0.33%  "

# Row 6 - XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.600"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.35%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.26%  "

# Row 8 - Solana
$ws.Range("D8").Value = "41.13"
$ws.Range("E8").Value = "  +13.30%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.292"
$ws.Range("E9").Value = "  -0.11%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.0667"
$ws.Range("E10").Value = "  -1.20%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.0998"
$ws.Range("E11").Value = "  +3.61%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.061.17"
$ws.Range("E12").Value = "  +0.24%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.792.16"
$ws.Range("E13").Value = "  -0.27%  "

# Row 14 - Chainlink
$ws.Range("D14").Value = "10.95"
$ws.Range("E14").Value = "  -2.95%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "34.485.85"
$ws.Range("E15").Value = "  +0.56%  "

# Row 16 - Polygon
$ws.Range("D16").Value = "0.628"
$ws.Range("E16").Value = "  -0.73%  "

# Row 17 - Polkadot
$ws.Range("E17").Value = "  -0.25%  "

# Row 18 - Litecoin
$ws.Range("D18").Value = "67.21"
$ws.Range("E18").Value = "  -2.02%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "240.07"
$ws.Range("E19").Value = "  -1.70%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0767"
$ws.Range("E20").Value = "  -0.51%  "

# Row 21 - Avalanche
$ws.Range("D21").Value = "11.12"
$ws.Range("E21").Value = "  -2.04%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  +0.28%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "4.22"
$ws.Range("E23").Value = "  +3.71%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -1.54%  "

# Row 25 - Monero
$ws.Range("D25").Value = "172.06"

# Row 26 - Cosmos
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.51%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "17.39"
$ws.Range("E27").Value = "  +0.59%  "

# Row 28 - Stellar
$ws.Range("D28").Value = "0.121"
$ws.Range("E28").Value = "  +0.22%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.35%  "

# Row 30 - Filecoin
$ws.Range("D30").Value = "3.78"
$ws.Range("E30").Value = "  +0.09%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -0.49%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  -0.78%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "0.0512"
$ws.Range("E33").Value = "  +0.04%  "

# Row 34 - LidoDAOToken
$ws.Range("D34").Value = "1.79"
$ws.Range("E34").Value = "  +1.24%  "

# Row 35 - was Maker, now ImmutableX
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "0.647"
$ws.Range("E35").Value = "  +0.29%  "

# Row 36 - was ImmutableX, now Maker
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "1.322.84"
$ws.Range("E36").Value = "  -2.77%  "

# Row 37 - TrustWalletToken
$ws.Range("E37").Value = "  +0.67%  "

# Row 38 - Aave
$ws.Range("D38").Value = "86.11"
$ws.Range("E38").Value = "  +6.85%  "

# Row 39 - was RenderToken, now VeChain
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.0188"
$ws.Range("E39").Value = "  +1.25%  "

# Row 40 - was InjectiveProtocol, now RenderToken
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "2.34"
$ws.Range("E40").Value = "  -0.49%  "

# Row 41 - was VeChain, now InjectiveProtocol
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.59%  "

# Row 42 - WEMIXToken
$ws.Range("E42").Value = "  +5.94%  "

# Row 43 - HuobiToken
$ws.Range("E43").Value = "  +0.70%  "

# Row 44 - MXToken
$ws.Range("E44").Value = "  +0.20%  "

# Row 45 - ARBITRUM
$ws.Range("E45").Value = "  +0.08%  "

# Row 46 - Kaspa
$ws.Range("D46").Value = "0.0519"
$ws.Range("E46").Value = "  +4.64%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "1.961.87"
$ws.Range("E47").Value = "  +0.25%  "

# Row 48 - FraxShare
$ws.Range("E48").Value = "  +1.15%  "

# Row 49 - PaxDollar
$ws.Range("E49").Value = "  +0.24%  "

# Row 50 - Quant
$ws.Range("D50").Value = "100.47"
$ws.Range("E50").Value = "  -1.33%  "

# Row 51 - Cronos
$ws.Range("E51").Value = "  +1.00%  "
